$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MGIC")

# Cells that become "NA" (text) instead of numeric values
$ws.Range("J21").Value = "NA"
$ws.Range("J83").Value = "NA"
$ws.Range("J94").Value = "NA"
$ws.Range("J100").Value = "NA"
$ws.Range("J101").Value = "NA"

# Row 91 (Capital Expenditures) updated values
$ws.Range("D91").Value = -1400
$ws.Range("E91").Value = -800
$ws.Range("F91").Value = -1100
$ws.Range("G91").Value = -1000
$ws.Range("H91").Value = -500
$ws.Range("I91").Value = -500
$ws.Range("J91").Value = -500
